# Weekly update: a new price record (week of 2023-10-19) is prepended to the
# "Hortaliza, Feria Lagunitas de Puerto Montt - Zanahoria" data set.
# This shifts every existing data row (571..659) down by one (572..660) and
# inserts the new record into the freed-up row 571.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 571..659 down to 572..660, freeing row 571 for the new record.
$ws.Rows.Item(571).Insert()

# Populate the newly inserted row with the latest weekly observation.
$ws.Cells.Item(571, 1).Value  = 4
$ws.Cells.Item(571, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(571, 3).Value  = "Los Lagos"
$ws.Cells.Item(571, 4).Value  = 45218
$ws.Cells.Item(571, 5).Value  = 10
$ws.Cells.Item(571, 6).Value  = 100114013
$ws.Cells.Item(571, 7).Value  = "Zanahoria"
$ws.Cells.Item(571, 8).Value  = "Sin especificar"
$ws.Cells.Item(571, 9).Value  = "Primera"
$ws.Cells.Item(571, 10).Value = 300
$ws.Cells.Item(571, 11).Value = 8500
$ws.Cells.Item(571, 12).Value = 8500
$ws.Cells.Item(571, 13).Value = 8500
$ws.Cells.Item(571, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(571, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(571, 16).Value = 425
$ws.Cells.Item(571, 17).Value = 20
$ws.Cells.Item(571, 18).Value = "Hortaliza"
